$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Scratch sheet + helper cell used to force numeric-looking strings to be
# stored as real text (matching the source data) without leaving stray
# formatting behind on the destination cells. Create it FIRST so that no
# further sheet insert/remove operations happen while other sheet
# references are in use.
$scratch = $wb.Worksheets.Add()
$scratch.Name = "__scratch__"
$helper = $scratch.Range("A1")

$total = $wb.Worksheets.Item("总计")

# --- 1. Update the "总计" (Total) sheet: insert a new row for 2022-Q1 ---
$total.Rows("2:2").Insert()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 15
$total.Range("D2").Value = 18.12
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

# --- 2. Create the new "2022-Q1" sheet (copy formatting/shape from "2021-Q4") ---
$sample = $wb.Worksheets.Item("2021-Q4")
$sample.Copy($total)
$new = $wb.Worksheets.Item("2021-Q4 (2)")
$new.Name = "2022-Q1"
$new.Cells.ClearContents()
$new.Rows("17:24").Delete()

# --- Header row ---
$new.Range("B1").Value = "基金代码"
$new.Range("C1").Value = "基金名称"
$new.Range("D1").Value = "基金规模"
$new.Range("E1").Value = "股票总仓位"
$new.Range("F1").Value = "仓位占比"
$new.Range("G1").Value = "持有市值(亿元)"
$new.Range("H1").Value = "仓位排名"

# --- Data rows ---
$new.Range("A2").Value = 0
$helper.Value = "'005669"
$helper.Copy()
$new.Range("B2").PasteSpecial(-4163)
$new.Range("C2").Value = "前海开源公用事业行业股票"
$helper.Value = "'258.16"
$helper.Copy()
$new.Range("D2").PasteSpecial(-4163)
$helper.Value = "'94.53"
$helper.Copy()
$new.Range("E2").PasteSpecial(-4163)
$helper.Value = "'5.78"
$helper.Copy()
$new.Range("F2").PasteSpecial(-4163)
$helper.Value = "'14.9216"
$helper.Copy()
$new.Range("G2").PasteSpecial(-4163)
$new.Range("H2").Value = 4

$new.Range("A3").Value = 1
$helper.Value = "'009007"
$helper.Copy()
$new.Range("B3").PasteSpecial(-4163)
$new.Range("C3").Value = "兴全沪港深两年持有期混合"
$helper.Value = "'24.80"
$helper.Copy()
$new.Range("D3").PasteSpecial(-4163)
$helper.Value = "'90.13"
$helper.Copy()
$new.Range("E3").PasteSpecial(-4163)
$helper.Value = "'3.64"
$helper.Copy()
$new.Range("F3").PasteSpecial(-4163)
$helper.Value = "'0.9027"
$helper.Copy()
$new.Range("G3").PasteSpecial(-4163)
$new.Range("H3").Value = 2

$new.Range("A4").Value = 2
$helper.Value = "'009630"
$helper.Copy()
$new.Range("B4").PasteSpecial(-4163)
$new.Range("C4").Value = "浦银安盛ESG责任投资混合A"
$helper.Value = "'15.61"
$helper.Copy()
$new.Range("D4").PasteSpecial(-4163)
$helper.Value = "'80.10"
$helper.Copy()
$new.Range("E4").PasteSpecial(-4163)
$helper.Value = "'4.85"
$helper.Copy()
$new.Range("F4").PasteSpecial(-4163)
$helper.Value = "'0.7571"
$helper.Copy()
$new.Range("G4").PasteSpecial(-4163)
$new.Range("H4").Value = 7

$new.Range("A5").Value = 3
$helper.Value = "'012073"
$helper.Copy()
$new.Range("B5").PasteSpecial(-4163)
$new.Range("C5").Value = "华安均衡优选混合A"
$helper.Value = "'8.33"
$helper.Copy()
$new.Range("D5").PasteSpecial(-4163)
$helper.Value = "'89.23"
$helper.Copy()
$new.Range("E5").PasteSpecial(-4163)
$helper.Value = "'5.22"
$helper.Copy()
$new.Range("F5").PasteSpecial(-4163)
$helper.Value = "'0.4348"
$helper.Copy()
$new.Range("G5").PasteSpecial(-4163)
$new.Range("H5").Value = 1

$new.Range("A6").Value = 4
$helper.Value = "'009631"
$helper.Copy()
$new.Range("B6").PasteSpecial(-4163)
$new.Range("C6").Value = "浦银安盛ESG责任投资混合C"
$helper.Value = "'5.74"
$helper.Copy()
$new.Range("D6").PasteSpecial(-4163)
$helper.Value = "'80.10"
$helper.Copy()
$new.Range("E6").PasteSpecial(-4163)
$helper.Value = "'4.85"
$helper.Copy()
$new.Range("F6").PasteSpecial(-4163)
$helper.Value = "'0.2784"
$helper.Copy()
$new.Range("G6").PasteSpecial(-4163)
$new.Range("H6").Value = 7

$new.Range("A7").Value = 5
$helper.Value = "'010994"
$helper.Copy()
$new.Range("B7").PasteSpecial(-4163)
$new.Range("C7").Value = "博时创新经济混合A"
$helper.Value = "'3.89"
$helper.Copy()
$new.Range("D7").PasteSpecial(-4163)
$helper.Value = "'89.89"
$helper.Copy()
$new.Range("E7").PasteSpecial(-4163)
$helper.Value = "'6.55"
$helper.Copy()
$new.Range("F7").PasteSpecial(-4163)
$helper.Value = "'0.2548"
$helper.Copy()
$new.Range("G7").PasteSpecial(-4163)
$new.Range("H7").Value = 2

$new.Range("A8").Value = 6
$helper.Value = "'001581"
$helper.Copy()
$new.Range("B8").PasteSpecial(-4163)
$new.Range("C8").Value = "华安沪港深通精选灵活配置混合"
$helper.Value = "'4.92"
$helper.Copy()
$new.Range("D8").PasteSpecial(-4163)
$helper.Value = "'92.91"
$helper.Copy()
$new.Range("E8").PasteSpecial(-4163)
$helper.Value = "'4.86"
$helper.Copy()
$new.Range("F8").PasteSpecial(-4163)
$helper.Value = "'0.2391"
$helper.Copy()
$new.Range("G8").PasteSpecial(-4163)
$new.Range("H8").Value = 1

$new.Range("A9").Value = 7
$helper.Value = "'040018"
$helper.Copy()
$new.Range("B9").PasteSpecial(-4163)
$new.Range("C9").Value = "华安香港精选股票(QDII)"
$helper.Value = "'5.47"
$helper.Copy()
$new.Range("D9").PasteSpecial(-4163)
$helper.Value = "'88.46"
$helper.Copy()
$new.Range("E9").PasteSpecial(-4163)
$helper.Value = "'4.09"
$helper.Copy()
$new.Range("F9").PasteSpecial(-4163)
$helper.Value = "'0.2237"
$helper.Copy()
$new.Range("G9").PasteSpecial(-4163)
$new.Range("H9").Value = 2

$new.Range("A10").Value = 8
$helper.Value = "'009017"
$helper.Copy()
$new.Range("B10").PasteSpecial(-4163)
$new.Range("C10").Value = "银华港股通精选股票"
$helper.Value = "'0.91"
$helper.Copy()
$new.Range("D10").PasteSpecial(-4163)
$helper.Value = "'86.12"
$helper.Copy()
$new.Range("E10").PasteSpecial(-4163)
$helper.Value = "'5.69"
$helper.Copy()
$new.Range("F10").PasteSpecial(-4163)
$helper.Value = "'0.0518"
$helper.Copy()
$new.Range("G10").PasteSpecial(-4163)
$new.Range("H10").Value = 3

$new.Range("A11").Value = 9
$helper.Value = "'010995"
$helper.Copy()
$new.Range("B11").PasteSpecial(-4163)
$new.Range("C11").Value = "博时创新经济混合C"
$helper.Value = "'0.35"
$helper.Copy()
$new.Range("D11").PasteSpecial(-4163)
$helper.Value = "'89.89"
$helper.Copy()
$new.Range("E11").PasteSpecial(-4163)
$helper.Value = "'6.55"
$helper.Copy()
$new.Range("F11").PasteSpecial(-4163)
$helper.Value = "'0.0229"
$helper.Copy()
$new.Range("G11").PasteSpecial(-4163)
$new.Range("H11").Value = 2

$new.Range("A12").Value = 10
$helper.Value = "'012074"
$helper.Copy()
$new.Range("B12").PasteSpecial(-4163)
$new.Range("C12").Value = "华安均衡优选混合C"
$helper.Value = "'0.25"
$helper.Copy()
$new.Range("D12").PasteSpecial(-4163)
$helper.Value = "'89.23"
$helper.Copy()
$new.Range("E12").PasteSpecial(-4163)
$helper.Value = "'5.22"
$helper.Copy()
$new.Range("F12").PasteSpecial(-4163)
$helper.Value = "'0.0130"
$helper.Copy()
$new.Range("G12").PasteSpecial(-4163)
$new.Range("H12").Value = 1

$new.Range("A13").Value = 11
$helper.Value = "'040021"
$helper.Copy()
$new.Range("B13").PasteSpecial(-4163)
$new.Range("C13").Value = "华安大中华升级股票(QDII)"
$helper.Value = "'0.26"
$helper.Copy()
$new.Range("D13").PasteSpecial(-4163)
$helper.Value = "'87.37"
$helper.Copy()
$new.Range("E13").PasteSpecial(-4163)
$helper.Value = "'3.93"
$helper.Copy()
$new.Range("F13").PasteSpecial(-4163)
$helper.Value = "'0.0102"
$helper.Copy()
$new.Range("G13").PasteSpecial(-4163)
$new.Range("H13").Value = 3

$new.Range("A14").Value = 12
$helper.Value = "'501303"
$helper.Copy()
$new.Range("B14").PasteSpecial(-4163)
$new.Range("C14").Value = "广发港股通恒生综合中型股指数(LOF)A"
$helper.Value = "'0.34"
$helper.Copy()
$new.Range("D14").PasteSpecial(-4163)
$helper.Value = "'92.39"
$helper.Copy()
$new.Range("E14").PasteSpecial(-4163)
$helper.Value = "'1.99"
$helper.Copy()
$new.Range("F14").PasteSpecial(-4163)
$helper.Value = "'0.0068"
$helper.Copy()
$new.Range("G14").PasteSpecial(-4163)
$new.Range("H14").Value = 3

$new.Range("A15").Value = 13
$helper.Value = "'004996"
$helper.Copy()
$new.Range("B15").PasteSpecial(-4163)
$new.Range("C15").Value = "广发港股通恒生综合中型股指数(LOF)C"
$helper.Value = "'0.11"
$helper.Copy()
$new.Range("D15").PasteSpecial(-4163)
$helper.Value = "'92.39"
$helper.Copy()
$new.Range("E15").PasteSpecial(-4163)
$helper.Value = "'1.99"
$helper.Copy()
$new.Range("F15").PasteSpecial(-4163)
$helper.Value = "'0.0022"
$helper.Copy()
$new.Range("G15").PasteSpecial(-4163)
$new.Range("H15").Value = 3

$new.Range("A16").Value = 14
$helper.Value = "'160922"
$helper.Copy()
$new.Range("B16").PasteSpecial(-4163)
$new.Range("C16").Value = "大成恒生综合中小型股指数(QDII-LOF)A"
$helper.Value = "'0.10"
$helper.Copy()
$new.Range("D16").PasteSpecial(-4163)
$helper.Value = "'92.44"
$helper.Copy()
$new.Range("E16").PasteSpecial(-4163)
$helper.Value = "'1.60"
$helper.Copy()
$new.Range("F16").PasteSpecial(-4163)
$helper.Value = "'0.0016"
$helper.Copy()
$new.Range("G16").PasteSpecial(-4163)
$new.Range("H16").Value = 3

$scratch.Delete()

$wb.Worksheets.Item(1).Activate()
